$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 195
$ws.Range("S195").Value = "`"she under the influence of black magic (?) developed a psychic illness and ingested a knife`""
$ws.Range("AH195").Value = "`"There were omental adhesions all over the abdominal cavity extending up to pelvis. Duodenum, liver and hepatic flexure of the colon were densely adhered with each other, forming a sort of cocoon. Hub of the knife could be palpated at the cocooned site in the region of hepatic flexure.`", `"During adhenolysis there was a rent in the duodenum just at the hub of the knife`""
$ws.Range("B195").Value = "626-001"
$ws.Range("Z195").Value = "`"17.78cm knife`", `"in situ 7 years`""
$ws.Range("Z195").Font.Color = 0
$ws.Range("A195").Value = 626
$ws.Range("C195").Value = 32
$ws.Range("D195").Value = "Female"
$ws.Range("E195").Value = "Y"
$ws.Range("F195").Value = "N"
$ws.Range("G195").Value = "N"
$ws.Range("H195").Value = "UK"
$ws.Range("I195").Value = "UK"
$ws.Range("J195").Value = "N"
$ws.Range("K195").Value = "N"
$ws.Range("L195").Value = "N"
$ws.Range("M195").Value = "N"
$ws.Range("N195").Value = "N"
$ws.Range("O195").Value = "Y"
$ws.Range("P195").Value = "N"
$ws.Range("Q195").Value = "N"
$ws.Range("R195").Value = "N"
$ws.Range("T195").Value = "N"
$ws.Range("U195").Value = "N"
$ws.Range("V195").Value = "Y"
$ws.Range("W195").Value = "Y"
$ws.Range("X195").Value = "Y"
$ws.Range("Y195").Value = "N"
$ws.Range("AA195").Value = "Y"
$ws.Range("AB195").Value = "Y"
$ws.Range("AC195").Value = "N"
$ws.Range("AD195").Value = "Y"
$ws.Range("AE195").Value = "Y"
$ws.Range("AF195").Value = "N"
$ws.Range("AG195").Value = "Y"

# Row 196
$ws.Range("B196").Value = "637-001"
$ws.Range("S196").Value = "`"In the context of schizophrenia, the development of pica has various potential causes: • Long standing malnutrition related to an underlying chronic psychotic illness • Psychotropic induced compulsive eating behavior of inedible substances[4] • Hematopoietic suppression resulting from a. Chronic psychotic illness b. Chronic use of psychotropics. • A feature of disorganized behavior. Pica is believed to be an obsessive‑compulsi`""
$ws.Range("Z196").Value = "`"she was eating clay and pieces of brick for approximately 2 months prior to hospitalization`""
$ws.Range("Z196").Font.Color = 0
$ws.Range("A196").Value = 637
$ws.Range("C196").Value = 22
$ws.Range("D196").Value = "Female"
$ws.Range("E196").Value = "Y"
$ws.Range("F196").Value = "N"
$ws.Range("G196").Value = "N"
$ws.Range("H196").Value = "UK"
$ws.Range("I196").Value = "UK"
$ws.Range("J196").Value = "Y"
$ws.Range("K196").Value = "N"
$ws.Range("L196").Value = "N"
$ws.Range("M196").Value = "N"
$ws.Range("N196").Value = "N"
$ws.Range("O196").Value = "Y"
$ws.Range("P196").Value = "N"
$ws.Range("Q196").Value = "N"
$ws.Range("R196").Value = "N"
$ws.Range("T196").Value = "N"
$ws.Range("U196").Value = "N"
$ws.Range("V196").Value = "N"
$ws.Range("W196").Value = "N"
$ws.Range("X196").Value = "N"
$ws.Range("Y196").Value = "Y"
$ws.Range("AA196").Value = "N"
$ws.Range("AB196").Value = "N"
$ws.Range("AC196").Value = "N"
$ws.Range("AD196").Value = "N"
$ws.Range("AE196").Value = "N"
$ws.Range("AF196").Value = "N"
$ws.Range("AG196").Value = "N"

# Row 197
$ws.Range("S197").Value = "Background of paranoid schizophrenia, `"been swallowing metallic items “out of boredom” for the previous 12 months.`""
$ws.Range("B197").Value = "640-001"
$ws.Range("Z197").Value = "`"50 metal objects`""
$ws.Range("Z197").Font.Color = 0
$ws.Range("A197").Value = 640
$ws.Range("C197").Value = 24
$ws.Range("D197").Value = "Male"
$ws.Range("E197").Value = "Y"
$ws.Range("F197").Value = "N"
$ws.Range("G197").Value = "N"
$ws.Range("H197").Value = "UK"
$ws.Range("I197").Value = "UK"
$ws.Range("J197").Value = "Y"
$ws.Range("K197").Value = "N"
$ws.Range("L197").Value = "N"
$ws.Range("M197").Value = "N"
$ws.Range("N197").Value = "N"
$ws.Range("O197").Value = "Y"
$ws.Range("P197").Value = "Y"
$ws.Range("Q197").Value = "N"
$ws.Range("R197").Value = "N"
$ws.Range("T197").Value = "N"
$ws.Range("U197").Value = "N"
$ws.Range("V197").Value = "Y"
$ws.Range("W197").Value = "Y"
$ws.Range("X197").Value = "Y"
$ws.Range("Y197").Value = "Y"
$ws.Range("AA197").Value = "N"
$ws.Range("AB197").Value = "Y"
$ws.Range("AC197").Value = "N"
$ws.Range("AD197").Value = "Y"
$ws.Range("AE197").Value = "Y"
$ws.Range("AF197").Value = "N"
$ws.Range("AG197").Value = "N"

# Row 198
$ws.Range("S198").Value = "`"the patient’s swallowing of the lighters was associated with impulsive behavior due to mental retardation`""
$ws.Range("B198").Value = "643-001"
$ws.Range("Z198").Value = "`"12 lighters`""
$ws.Range("Z198").Font.Color = 0
$ws.Range("A198").Value = 643
$ws.Range("C198").Value = 21
$ws.Range("D198").Value = "Male"
$ws.Range("E198").Value = "Y"
$ws.Range("F198").Value = "N"
$ws.Range("G198").Value = "UK"
$ws.Range("H198").Value = "UK"
$ws.Range("I198").Value = "UK"
$ws.Range("J198").Value = "Y"
$ws.Range("K198").Value = "Y"
$ws.Range("L198").Value = "N"
$ws.Range("M198").Value = "N"
$ws.Range("N198").Value = "N"
$ws.Range("O198").Value = "Y"
$ws.Range("P198").Value = "N"
$ws.Range("Q198").Value = "N"
$ws.Range("R198").Value = "N"
$ws.Range("T198").Value = "N"
$ws.Range("U198").Value = "N"
$ws.Range("V198").Value = "Y"
$ws.Range("W198").Value = "Y"
$ws.Range("X198").Value = "N"
$ws.Range("Y198").Value = "Y"
$ws.Range("AA198").Value = "Y"
$ws.Range("AB198").Value = "Y"
$ws.Range("AC198").Value = "N"
$ws.Range("AD198").Value = "Y"
$ws.Range("AE198").Value = "Y"
$ws.Range("AF198").Value = "Y"
$ws.Range("AG198").Value = "N"

# Row 199
$ws.Range("S199").Value = "`"10-year-old male with autism`""
$ws.Range("Z199").Value = "`"31 spherical magnets`""
$ws.Range("Z199").Font.Color = 0
$ws.Range("AH199").Value = "`"Three areas of full thickness erosions were noted as well as two areas of partial erosion, but no segments of bowel required resection`""
$ws.Range("B199").Value = "644-001"
$ws.Range("A199").Value = 644
$ws.Range("C199").Value = 10
$ws.Range("D199").Value = "Male"
$ws.Range("E199").Value = "Y"
$ws.Range("F199").Value = "N"
$ws.Range("G199").Value = "UK"
$ws.Range("H199").Value = "UK"
$ws.Range("I199").Value = "N"
$ws.Range("J199").Value = "Y"
$ws.Range("K199").Value = "Y"
$ws.Range("L199").Value = "UK"
$ws.Range("M199").Value = "N"
$ws.Range("N199").Value = "N"
$ws.Range("O199").Value = "N"
$ws.Range("P199").Value = "N"
$ws.Range("Q199").Value = "Y"
$ws.Range("R199").Value = "N"
$ws.Range("T199").Value = "N"
$ws.Range("U199").Value = "Y"
$ws.Range("V199").Value = "N"
$ws.Range("W199").Value = "N"
$ws.Range("X199").Value = "N"
$ws.Range("Y199").Value = "Y"
$ws.Range("AA199").Value = "Y"
$ws.Range("AB199").Value = "Y"
$ws.Range("AC199").Value = "N"
$ws.Range("AD199").Value = "Y"
$ws.Range("AE199").Value = "Y"
$ws.Range("AF199").Value = "N"
$ws.Range("AG199").Value = "N"

# Row 200
$ws.Range("B200").Value = "651-001"
$ws.Range("Z200").Value = "`"15 plastic spoons`""
$ws.Range("Z200").Font.Color = 0
$ws.Range("S200").Value = "`"inmate with a history of schizophrenia`""
$ws.Range("AH200").Value = "`"ulceration`""
$ws.Range("A200").Value = 651
$ws.Range("C200").Value = 34
$ws.Range("D200").Value = "Male"
$ws.Range("E200").Value = "Y"
$ws.Range("F200").Value = "Y"
$ws.Range("G200").Value = "N"
$ws.Range("H200").Value = "UK"
$ws.Range("I200").Value = "UK"
$ws.Range("J200").Value = "Y"
$ws.Range("K200").Value = "N"
$ws.Range("L200").Value = "UK"
$ws.Range("M200").Value = "UK"
$ws.Range("N200").Value = "UK"
$ws.Range("O200").Value = "UK"
$ws.Range("P200").Value = "UK"
$ws.Range("Q200").Value = "Y"
$ws.Range("T200").Value = "N"
$ws.Range("U200").Value = "N"
$ws.Range("V200").Value = "Y"
$ws.Range("W200").Value = "Y"
$ws.Range("X200").Value = "N"
$ws.Range("Y200").Value = "Y"
$ws.Range("AA200").Value = "Y"
$ws.Range("AB200").Value = "N"
$ws.Range("AC200").Value = "N"
$ws.Range("AD200").Value = "Y"
$ws.Range("AE200").Value = "N"
$ws.Range("AF200").Value = "N"
$ws.Range("AG200").Value = "Y"

# Final selection (matches the author's last-saved cursor position)
$ws.Range("AF200").Select() | Out-Null
